$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B2").Value = 26
$ws.Range("B3").Value = 1857142.857142857
$ws.Range("B13").Value = 2785714.285714286
$ws.Range("B23").Value = 2785714.285714286
$ws.Range("B31").Value = 1857142.857142857
$ws.Range("B32").Value = 2785714.285714286
$ws.Range("B33").Value = 2785714.285714286
$ws.Range("B34").Value = 7428571.428571429
